$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.315.84'
$ws.Range('E2').Value = '  -2.78%  '
$ws.Range('D3').Value = '3.683.60'
$ws.Range('E3').Value = '  -3.43%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = "'684.34"
$ws.Range('E5').Value = '  -2.63%  '
$ws.Range('D6').Value = "'159.98"
$ws.Range('E6').Value = '  -6.53%  '
$ws.Range('D7').Value = '3.680.86'
$ws.Range('E7').Value = '  -3.50%  '
$ws.Range('E8').Value = '  -0.12%  '
$ws.Range('E9').Value = '  -6.29%  '
$ws.Range('D10').Value = "'0.145"
$ws.Range('E10').Value = '  -9.11%  '
$ws.Range('D11').Value = "'7.19"
$ws.Range('E11').Value = '  -3.65%  '
$ws.Range('E12').Value = '  -10.84%  '
$ws.Range('E13').Value = '  -7.36%  '
$ws.Range('D14').Value = '4.304.23'
$ws.Range('E14').Value = '  -3.44%  '
$ws.Range('D15').Value = "'32.32"
$ws.Range('E15').Value = '  -11.99%  '
$ws.Range('D16').Value = '3.691.40'
$ws.Range('E16').Value = '  -3.46%  '
$ws.Range('D17').Value = '69.337.53'
$ws.Range('E17').Value = '  -2.86%  '
$ws.Range('E18').Value = '  -1.42%  '
$ws.Range('D19').Value = "'15.84"
$ws.Range('E19').Value = '  -9.70%  '
$ws.Range('D20').Value = "'6.42"
$ws.Range('E20').Value = '  -11.32%  '
$ws.Range('D21').Value = "'472.32"
$ws.Range('E21').Value = '  -7.56%  '
$ws.Range('E22').Value = '  -5.65%  '
$ws.Range('D23').Value = "'0.645"
$ws.Range('E23').Value = '  -9.95%  '
$ws.Range('D24').Value = "'79.52"
$ws.Range('E24').Value = '  -5.10%  '
$ws.Range('D25').Value = '3.829.27'
$ws.Range('E25').Value = '  -3.32%  '
$ws.Range('E26').Value = '  -0.13%  '
$ws.Range('E27').Value = '  -12.31%  '
$ws.Range('E28').Value = '  -14.97%  '
$ws.Range('E29').Value = '  -11.55%  '
$ws.Range('E30').Value = '  -10.74%  '
$ws.Range('E31').Value = '  -13.37%  '
$ws.Range('B32').Value = 'NEARProtocol'
$ws.Range('C32').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D32').Value = "'6.63"
$ws.Range('E32').Value = '  -10.02%  '
$ws.Range('B33').Value = 'ImmutableX'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D33').Value = "'2.02"
$ws.Range('E33').Value = '  -10.96%  '
$ws.Range('E34').Value = '  +0.13%  '
$ws.Range('D35').Value = "'26.61"
$ws.Range('E35').Value = '  -8.89%  '
$ws.Range('D36').Value = "'0.158"
$ws.Range('E36').Value = '  -7.77%  '
$ws.Range('D37').Value = "'8.15"
$ws.Range('E37').Value = '  -12.43%  '
$ws.Range('E38').Value = '  -9.37%  '
$ws.Range('D39').Value = "'2.26"
$ws.Range('E39').Value = '  -5.87%  '
$ws.Range('D41').Value = "'0.0898"
$ws.Range('E41').Value = '  -11.25%  '
$ws.Range('E42').Value = '  +0.09%  '
$ws.Range('D43').Value = "'0.939"
$ws.Range('E43').Value = '  -7.21%  '
$ws.Range('D44').Value = "'164.96"
$ws.Range('E44').Value = '  -0.95%  '
$ws.Range('D45').Value = "'47.86"
$ws.Range('E45').Value = '  -4.39%  '
$ws.Range('D46').Value = "'2.70"
$ws.Range('E46').Value = '  -16.17%  '
$ws.Range('D47').Value = "'1.30"
$ws.Range('E47').Value = '  -6.83%  '
$ws.Range('E48').Value = '  -5.32%  '
$ws.Range('B49').Value = 'InjectiveProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D49').Value = "'28.14"
$ws.Range('E49').Value = '  -8.33%  '
$ws.Range('B50').Value = 'FLOKI'
$ws.Range('C50').Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range('D50').Value = "'0.000272"
$ws.Range('E50').Value = '  -10.51%  '
$ws.Range('D51').Value = "'7.83"
$ws.Range('E51').Value = '  -9.59%  '
